# Contest 4 RR vs PBKS
# - Rename two team/predictor headers on Sheet1
# - Fill in prediction scores for matches 2-4 (rows 14-16)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename predictor/team headers in row 11 (merged header cells)
$ws.Range("M11").Value = "CheemsRajah"
$ws.Range("AB11").Value = "Rag Nat0112"

# Match 2 (row 14) - CSK vs DC
$ws.Range("E14").Value = 50
$ws.Range("H14").Value = 20
$ws.Range("K14").Value = 40
$ws.Range("N14").Value = 60
$ws.Range("Q14").Value = 0
$ws.Range("T14").Value = 70
$ws.Range("W14").Value = 100
$ws.Range("Z14").Value = 30
$ws.Range("AC14").Value = 80

# Match 3 (row 15) - SRH vs KKR
$ws.Range("E15").Value = 80
$ws.Range("H15").Value = 40
$ws.Range("K15").Value = 70
$ws.Range("N15").Value = 20
$ws.Range("Q15").Value = 100
$ws.Range("T15").Value = 0
$ws.Range("W15").Value = 50
$ws.Range("Z15").Value = 60
$ws.Range("AC15").Value = 30

# Match 4 (row 16) - RR vs PBKS
$ws.Range("E16").Value = 30
$ws.Range("H16").Value = 20
$ws.Range("K16").Value = 0
$ws.Range("N16").Value = 70
$ws.Range("Q16").Value = 60
$ws.Range("T16").Value = 100
$ws.Range("W16").Value = 50
$ws.Range("Z16").Value = 80
$ws.Range("AC16").Value = 40
